$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("TipoOperacion") is renamed to "Tipo Operación", and its two
# data rows switch from the old "TipoOperacion" placeholder text to the
# actual operation type "SIMPLE".
$ws.Range("G1").Value = "Tipo Operación"
$ws.Range("G2").Value = "SIMPLE"
$ws.Range("G3").Value = "SIMPLE"

# Widen column G so the longer header/value fits (matches column J's width).
$ws.Range("G1").ColumnWidth = $ws.Range("J1").ColumnWidth

# "Numero Propuesta" values are refreshed with newer proposal numbers.
# Re-enter them as plain text (no quote-prefix) and reset to the default
# style so they no longer carry the old quotePrefix formatting.
$ws.Range("U2").Value = "'4899839"
$ws.Range("U2").Style = "Normal"
$ws.Range("U3").Value = "'4899840"
$ws.Range("U3").Style = "Normal"

# The "Resultado" cells keep their text but drop their explicit (redundant)
# style, falling back to the row's default formatting.
$ws.Range("V2").Style = "Normal"
$ws.Range("V3").Style = "Normal"

# Leave the selection on U6, matching where the editor's cursor ended up.
$ws.Range("U6").Select()
